$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 3 new header cells (UI Week 26/27/28) in AB1:AD1 ---
$ws.Range("AB1").Value = "UI Week 26"
$ws.Range("AC1").Value = "UI Week 27"
$ws.Range("AD1").Value = "UI Week 28"

# --- Per-county data: new weekly figures for weeks 26-28, and the
#     updated running Grand Total (old AA value + the three new weeks). ---
$data = @"
2,45,51,59,1594
3,81,63,68,2518
4,205,277,214,6374
5,66,53,68,1797
6,963,982,1237,35246
7,38,45,43,1310
8,48,58,71,1600
9,284,283,328,8158
10,197,232,274,7782
11,136,248,283,4283
12,291,362,419,10143
13,43,49,62,1917
14,2142,2616,2712,65485
15,523,610,677,12484
16,72,119,150,3705
17,114,129,138,3802
18,120,154,170,4724
19,360,427,504,13149
20,13,26,21,524
21,428,580,561,15024
22,29,44,55,1462
23,133,175,219,5300
24,172,198,194,4362
25,83,70,121,2936
26,93,174,119,4043
27,18,19,35,734
28,82,120,169,3147
29,350,333,419,10393
30,145,300,293,4584
31,495,542,637,17226
32,49,72,92,2431
33,767,595,543,16541
34,46,83,60,1742
35,112,104,123,2352
36,106,140,156,4489
37,520,495,641,13239
38,680,847,761,22439
39,299,419,416,7578
40,66,70,80,2659
41,20,27,23,732
42,4537,5270,6166,140529
43,210,389,383,6079
44,122,143,192,4569
45,122,153,207,4962
46,818,782,957,24438
47,315,427,503,10310
48,19,20,34,695
49,87,113,177,3454
50,124,143,196,4508
51,368,493,461,10734
52,86,102,119,2202
53,873,841,1022,26926
54,72,61,96,2130
55,857,902,880,23454
56,64,66,73,2805
57,504,376,662,12657
58,57,60,58,2293
59,145,237,203,5132
60,386,561,681,14671
61,287,409,459,8838
62,92,73,94,2453
63,107,110,146,4739
64,1403,1607,1699,78767
65,138,172,181,3444
66,50,64,76,2411
67,319,309,434,10335
68,63,71,71,2086
69,619,796,863,18900
70,1545,1642,1997,48960
71,304,362,572,7914
72,176,115,134,3214
73,803,834,865,22405
74,250,314,361,8948
75,26356,30208,34237,855970
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $r  = [int]$parts[0]
    $aa = [int]$parts[1]
    $ab = [int]$parts[2]
    $ac = [int]$parts[3]
    $ad = [int]$parts[4]
    $ws.Cells.Item($r, 27).Value = $aa
    $ws.Cells.Item($r, 28).Value = $ab
    $ws.Cells.Item($r, 29).Value = $ac
    $ws.Cells.Item($r, 30).Value = $ad
}

# --- Match the new column widths used for the extra weekly columns ---
$ws.Range("AB1:AD1").ColumnWidth = 10.5

# --- Restore the selected cell noted in the saved workbook state ---
[void]$ws.Range("Q57").Select()
